# Update the DASP-to-SWC ranking table (A2:C10) with the new ordering / counts
# as described in the commit: "Acrescentando a ligação da vulnerabilidade do
# SWC-112 para o DASP-2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: new values for columns A (rank), B (label), C (count)
$data = @(
    @(2,  1, "arithmetic",          25),
    @(3,  3, "reentrancy",          24),
    @(4,  6, "front_running",       18),
    @(5,  9, "Other",               18),
    @(6,  0, "access_control",      17),
    @(7,  4, "unchecked_low_calls", 11),
    @(8,  7, "time_manipulation",   7),
    @(9,  5, "bad_randomness",      4),
    @(10, 2, "denial_service",      2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
